$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H64").Value = 7504.6943
$ws.Range("I64").Value = 4166.8
$ws.Range("J64").Value = 8788.5
$ws.Range("K64").Value = 4166.8
$ws.Range("L64").Value = 8788.5
$ws.Range("M64").Value = -3918.8
$ws.Range("N64").Value = -9284.5
$ws.Range("H67").Value = 7504.6943
$ws.Range("I67").Value = 4166.8
$ws.Range("J67").Value = 8788.5
$ws.Range("K67").Value = 4166.8
$ws.Range("L67").Value = 8788.5
$ws.Range("M67").Value = -3308.8
$ws.Range("N67").Value = -10504.5
$ws.Range("H69").Value = 250006000
$ws.Range("J69").Value = 250006000
$ws.Range("L69").Value = 750018000
$ws.Range("N69").Value = -750019748
$ws.Range("H72").Value = 250006000
$ws.Range("J72").Value = 250006000
$ws.Range("L72").Value = 2250054000
$ws.Range("N72").Value = -2250062736
$ws.Range("H74").Value = 4800
$ws.Range("I74").Value = 4800
$ws.Range("K74").Value = 4800
$ws.Range("M74").Value = -3864
$ws.Range("H76").Value = 4998.8
$ws.Range("I76").Value = 4996.6665
$ws.Range("J76").Value = 5002
$ws.Range("K76").Value = 4996.6665
$ws.Range("L76").Value = 5002
$ws.Range("M76").Value = -4681.6665
$ws.Range("N76").Value = -5632
$ws.Range("H77").Value = 4800
$ws.Range("I77").Value = 4800
$ws.Range("K77").Value = 24000
$ws.Range("M77").Value = -19320
$ws.Range("H79").Value = 4998.8
$ws.Range("I79").Value = 4996.6665
$ws.Range("J79").Value = 5002
$ws.Range("K79").Value = 4996.6665
$ws.Range("L79").Value = 5002
$ws.Range("M79").Value = -3904.6665
$ws.Range("N79").Value = -7186
$ws.Range("H115").Value = 1911
$ws.Range("I115").Value = 729.5
$ws.Range("K115").Value = 2188.5
$ws.Range("M115").Value = -621.5
$ws.Range("H141").Value = 8725.889
$ws.Range("I141").Value = 8215.4
$ws.Range("J141").Value = 8922.23
$ws.Range("K141").Value = 24646.2
$ws.Range("L141").Value = 26766.69
$ws.Range("M141").Value = -19466.2
$ws.Range("N141").Value = -37126.69

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10626.637
$ws.Range("I32").Value = 12210.667
$ws.Range("J32").Value = 3498.5
$ws.Range("K32").Value = 12210.667
$ws.Range("L32").Value = 3498.5
$ws.Range("M32").Value = -11923.667
$ws.Range("N32").Value = -4072.5
$ws.Range("H61").Value = 2804.0833
$ws.Range("I61").Value = 2804.0833
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2804.0833
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2592.0833
$ws.Range("N61").ClearContents()
$ws.Range("H136").Value = 2804.0833
$ws.Range("I136").Value = 2804.0833
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8412.249899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5862.249899999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3141.8635
$ws.Range("I107").Value = 1258.2667
$ws.Range("J107").Value = 7178.143
$ws.Range("K107").Value = 1258.2667
$ws.Range("L107").Value = 7178.143
$ws.Range("M107").Value = 661.7333000000001
$ws.Range("N107").Value = -11018.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5460.2666
$ws.Range("I31").Value = 1809.7273
$ws.Range("J31").Value = 15499.25
$ws.Range("K31").Value = 1809.7273
$ws.Range("L31").Value = 15499.25
$ws.Range("M31").Value = -1514.7273
$ws.Range("N31").Value = -16089.25
$ws.Range("H34").Value = 5460.2666
$ws.Range("I34").Value = 1809.7273
$ws.Range("J34").Value = 15499.25
$ws.Range("K34").Value = 1809.7273
$ws.Range("L34").Value = 15499.25
$ws.Range("M34").Value = -1607.7273
$ws.Range("N34").Value = -15903.25
$ws.Range("H62").Value = 6138.8
$ws.Range("I62").Value = 6424.5
$ws.Range("J62").Value = 4996
$ws.Range("K62").Value = 6424.5
$ws.Range("L62").Value = 4996
$ws.Range("M62").Value = -5800.5
$ws.Range("N62").Value = -6244
$ws.Range("H65").Value = 6138.8
$ws.Range("I65").Value = 6424.5
$ws.Range("J65").Value = 4996
$ws.Range("K65").Value = 32122.5
$ws.Range("L65").Value = 24980
$ws.Range("M65").Value = -29002.5
$ws.Range("N65").Value = -31220
$ws.Range("H132").Value = 1893.2059
$ws.Range("I132").Value = 1900.1875
$ws.Range("J132").Value = 1781.5
$ws.Range("K132").Value = 5700.5625
$ws.Range("L132").Value = 5344.5
$ws.Range("M132").Value = -3170.5625
$ws.Range("N132").Value = -10404.5
$ws.Range("H134").Value = 2737.8096
$ws.Range("I134").Value = 2721.8462
$ws.Range("J134").Value = 2763.75
$ws.Range("K134").Value = 8165.5386
$ws.Range("L134").Value = 8291.25
$ws.Range("M134").Value = -5630.5386
$ws.Range("N134").Value = -13361.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4499
$ws.Range("J80").Value = 4499
$ws.Range("L80").Value = 13497
$ws.Range("N80").Value = -15369
$ws.Range("H83").Value = 4499
$ws.Range("J83").Value = 4499
$ws.Range("L83").Value = 40491
$ws.Range("N83").Value = -49851
$ws.Range("H87").Value = 5674.6665
$ws.Range("I87").Value = 5674.6665
$ws.Range("K87").Value = 17023.9995
$ws.Range("M87").Value = -15775.9995
$ws.Range("H90").Value = 5674.6665
$ws.Range("I90").Value = 5674.6665
$ws.Range("K90").Value = 51071.9985
$ws.Range("M90").Value = -44831.9985
$ws.Range("H113").Value = 1301.6111
$ws.Range("J113").Value = 1418.6
$ws.Range("L113").Value = 4255.799999999999
$ws.Range("N113").Value = -8595.8
$ws.Range("H132").Value = 2369.1538
$ws.Range("J132").Value = 2600
$ws.Range("L132").Value = 23400
$ws.Range("N132").Value = -28460
$ws.Range("H140").Value = 67968.07
$ws.Range("I140").Value = 72519
$ws.Range("J140").Value = 4255
$ws.Range("K140").Value = 217557
$ws.Range("L140").Value = 12765
$ws.Range("M140").Value = -212377
$ws.Range("N140").Value = -23125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 515000
$ws.Range("J52").Value = 515000
$ws.Range("L52").Value = 515000
$ws.Range("N52").Value = -515518
$ws.Range("H80").Value = 7170.9585
$ws.Range("J80").Value = 2771.2856
$ws.Range("N80").Value = -4767.2856
$ws.Range("H83").Value = 7170.9585
$ws.Range("J83").Value = 2771.2856
$ws.Range("L83").Value = 13856.428
$ws.Range("N83").Value = -23840.428
$ws.Range("H126").Value = 3998
$ws.Range("I126").Value = 3998
$ws.Range("K126").Value = 11994
$ws.Range("M126").Value = -9524

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4757.143
$ws.Range("I7").Value = 4900
$ws.Range("K7").Value = 4900
$ws.Range("M7").Value = -4788
$ws.Range("H68").Value = 12179.5
$ws.Range("I68").Value = 5001.3335
$ws.Range("K68").Value = 5001.3335
$ws.Range("M68").Value = -4252.3335
$ws.Range("H71").Value = 12179.5
$ws.Range("I71").Value = 5001.3335
$ws.Range("K71").Value = 25006.6675
$ws.Range("M71").Value = -21262.6675
$ws.Range("H126").Value = 4757.143
$ws.Range("I126").Value = 4900
$ws.Range("K126").Value = 14700
$ws.Range("M126").Value = -12230

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 14000
$ws.Range("J11").Value = 14000
$ws.Range("L11").Value = 14000
$ws.Range("N11").Value = -14284
$ws.Range("H126").Value = 1120
$ws.Range("I126").Value = 1120
$ws.Range("K126").Value = 3360
$ws.Range("M126").Value = -890
